# "Actualizacion de ph y de excel"
# Adds the new character "Milla" (row 2) to the Hoja1 sheet, stamps an
# (empty, underline-formatted) cell at C7, and updates the page setup /
# selection to match the saved workbook state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- New character row -----------------------------------------------
$ws.Range("A2").Value = "Rampo Doyle"
$ws.Range("B2").Value = "Milla"

# Puntos de vida, Proficiency bonus, Walking, Iniciativa, Clase de armadura,
# Strength..Charisma, Salvacion:STR..CHA, Passive checks, and all the skill
# scores (Acrobatics .. Survival) for the new character, columns C..AO.
$rowValues = @(
    33, 33, 3, 30, 0, 11, -1, 0, 1, 3, 0, 4, -1, 3, 0, 3, 1, 7, 10, 16,
    10, 0, 0, 6, -1, 7, 3, 0, 4, 6, 0, 3, 0, 4, 7, 3, 0, 3, 0
)
$col = 3
foreach ($v in $rowValues) {
    $ws.Cells.Item(2, $col).Value = $v
    $col++
}

# --- Empty, underlined placeholder cell -------------------------------
$ws.Range("C7").Font.Underline = $true

# --- Page setup ---------------------------------------------------------
$ws.PageSetup.PaperSize = 9    # xlPaperA4
$ws.PageSetup.Orientation = 1  # xlPortrait

# --- Selection / view state ---------------------------------------------
$ws.Range("C7").Select()
$excel.ActiveWindow.Left = 390
$excel.ActiveWindow.Top = 390
